# Append: 2025-12-22 12:51 JST
# Update the "取得日時" (acquisition timestamp) column on the "ランサーズ" sheet
# for all existing data rows (2-11) from the old timestamp to the new one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldTimestamp = "2025-12-22 12:38:44"
$newTimestamp = "2025-12-22 12:51:21"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
